# Updated cryptos list refresh: new Price / Volume(1h) figures pulled for
# each coin, plus two pairs of rows that swapped rank position (rows 27/28
# and rows 48/49 each exchanged their Coin/Link/Price/Volume data).
#
# Note: several Price values are purely-numeric-looking strings (e.g.
# "188.67") that must stay stored as TEXT (matching the sheet's existing
# convention of inline-string prices like "75.055.24"). Setting
# NumberFormat = "@" right before the assignment keeps Excel from silently
# coercing those into numeric cells.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Bitcoin ---
$ws.Range("D2").Value = "75.055.24"
$ws.Range("E2").Value = "  +3.18%  "

# --- Row 3: Ethereum ---
$ws.Range("D3").Value = "2.809.27"
$ws.Range("E3").Value = "  +9.08%  "

# --- Row 4: TetherUSD ---
$ws.Range("E4").Value = "  -0.04%  "

# --- Row 5: Solana ---
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "188.67"
$ws.Range("E5").Value = "  +4.61%  "

# --- Row 6: BNB ---
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "594.84"
$ws.Range("E6").Value = "  +3.20%  "

# --- Row 7: USDC ---
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  -0.04%  "

# --- Row 8: XRP ---
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.554"
$ws.Range("E8").Value = "  +5.17%  "

# --- Row 9: Dogecoin ---
$ws.Range("E9").Value = "  +1.59%  "

# --- Row 10: LidoStakedEther ---
$ws.Range("D10").Value = "2.806.93"
$ws.Range("E10").Value = "  +9.05%  "

# --- Row 11: TRON ---
$ws.Range("E11").Value = "  -0.42%  "

# --- Row 12: Cardano ---
$ws.Range("E12").Value = "  +4.43%  "

# --- Row 13: Toncoin ---
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.82"
$ws.Range("E13").Value = "  +3.11%  "

# --- Row 14: WrappedliquidstakedEther2.0 ---
$ws.Range("D14").Value = "3.330.89"
$ws.Range("E14").Value = "  +7.02%  "

# --- Row 15: WrappedBTC ---
$ws.Range("D15").Value = "74.896.97"
$ws.Range("E15").Value = "  +2.77%  "

# --- Row 16: ShibaInu ---
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000189"
$ws.Range("E16").Value = "  +3.07%  "

# --- Row 17: Avalanche ---
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "27.05"
$ws.Range("E17").Value = "  +5.82%  "

# --- Row 18: WrappedEther ---
$ws.Range("D18").Value = "2.819.61"
$ws.Range("E18").Value = "  +8.69%  "

# --- Row 19: Uniswap ---
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.08"
$ws.Range("E19").Value = "  +4.80%  "

# --- Row 20: Chainlink ---
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.28"
$ws.Range("E20").Value = "  +5.86%  "

# --- Row 21: BitcoinCash ---
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.69"
$ws.Range("E21").Value = "  +3.06%  "

# --- Row 22: SuiNetwork ---
$ws.Range("E22").Value = "  +3.95%  "

# --- Row 23: Polkadot ---
$ws.Range("E23").Value = "  +2.15%  "

# --- Row 24: LEO ---
$ws.Range("E24").Value = "  +0.10%  "

# --- Row 25: Dai ---
$ws.Range("E25").Value = "  -0.03%  "

# --- Row 26: Litecoin ---
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "70.90"
$ws.Range("E26").Value = "  +2.98%  "

# --- Rows 27 & 28 swapped: row 27 (was WrappedeETH) becomes NEARProtocol,
#     row 28 (was NEARProtocol) becomes WrappedeETH ---
$ws.Range("B27").Value = "NEARProtocol"
$ws.Range("C27").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.18"
$ws.Range("E27").Value = "  +3.33%  "

$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.968.14"
$ws.Range("E28").Value = "  +9.04%  "

# --- Row 29: Aptos ---
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "9.58"
$ws.Range("E29").Value = "  +5.97%  "

# --- Row 30: PEPE ---
$ws.Range("E30").Value = "  +13.88%  "

# --- Row 31: Binance-PegBSC-USD ---
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -1.27%  "

# --- Row 32: Bittensor ---
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "516.91"
$ws.Range("E32").Value = "  +5.50%  "

# --- Row 33: Fetch.AI ---
$ws.Range("E33").Value = "  +3.36%  "

# --- Row 34: InternetComputer(DFINITY) ---
$ws.Range("E34").Value = "  +2.95%  "

# --- Row 35: PancakeSwap ---
$ws.Range("E35").Value = "  +5.15%  "

# --- Row 36: FirstDigitalUSD ---
$ws.Range("E36").Value = "  -0.06%  "

# --- Row 37: Monero ---
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "163.19"
$ws.Range("E37").Value = "  +1.95%  "

# --- Row 38: EthereumClassic ---
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.09"
$ws.Range("E38").Value = "  +6.55%  "

# --- Row 39: Kaspa ---
$ws.Range("E39").Value = "  +0.94%  "

# --- Row 40: WhiteBITCoin ---
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "19.34"
$ws.Range("E40").Value = "  +0.54%  "

# --- Row 41: Aave ---
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "183.57"
$ws.Range("E41").Value = "  +19.54%  "

# --- Row 42: USDe ---
$ws.Range("E42").Value = "  +0.03%  "

# --- Row 43: RenderToken ---
$ws.Range("E43").Value = "  +6.54%  "

# --- Row 44: PolygonEcosystemToken ---
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.340"
$ws.Range("E44").Value = "  +7.07%  "

# --- Row 45: Stacks ---
$ws.Range("E45").Value = "  +2.86%  "

# --- Row 46: ImmutableX ---
$ws.Range("E46").Value = "  +6.61%  "

# --- Row 47: OKB ---
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "40.10"
$ws.Range("E47").Value = "  +4.27%  "

# --- Rows 48 & 49 swapped: row 48 (was Cronos) becomes dogwifhat,
#     row 49 (was dogwifhat) becomes Cronos ---
$ws.Range("B48").Value = "dogwifhat"
$ws.Range("C48").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.34"
$ws.Range("E48").Value = "  +3.04%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0861"
$ws.Range("E49").Value = "  +0.51%  "

# --- Row 50: ARBITRUM ---
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.568"
$ws.Range("E50").Value = "  +10.14%  "

# --- Row 51: Filecoin ---
$ws.Range("E51").Value = "  +5.06%  "
